$d = $word.ActiveDocument
$cr = [char]13

# ---------------------------------------------------------------------------
# 1. "Profesor: Rodrigo Núñez " + "Núñez" + "                       Grupo 1"
#    -> single run "Profesor: Rodrigo Núñez Núñez                       Grupo 1"
#    The runs on both sides of each <w:proofErr/> pair are consumed by the
#    match, so Word's merge naturally drops the now-orphaned proofErr tags.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Profesor: Rodrigo Núñez Núñez                       Grupo 1",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Profesor: Rodrigo Núñez Núñez                       Grupo 1", 2)

# ---------------------------------------------------------------------------
# 2. "Algoritmo Seleccionado: " + "Bubble" + " " + "Sort"
#    -> single run "Algoritmo Seleccionado: Bubble Sort"
#    Here the final run ("Sort") is the last child of its paragraph, so its
#    closing <w:proofErr w:type="spellEnd"/> can never sit "between" two runs
#    being merged in a single pass. Work around this by temporarily splitting
#    a new paragraph right after it, merging across that paragraph mark (which
#    relocates the dangling proofErr into the new, now-isolated paragraph),
#    and finally deleting that throw-away paragraph.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Algoritmo Seleccionado: Bubble",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Algoritmo Seleccionado: Bubble", 2)

$d.Content.Find.Execute(
    " Sort", $true, $false, $false, $false, $false, $true, 1, $false,
    " Sort", 2)

$algoritmoIndex = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text.TrimEnd($cr) -eq "Algoritmo Seleccionado: Bubble Sort") {
        $algoritmoIndex = $i
    }
}

$pAlgoritmo = $d.Paragraphs.Item($algoritmoIndex)
$pAlgoritmo.Range.InsertParagraphAfter()

$d.Content.Find.Execute(
    "Bubble Sort" + $cr, $true, $false, $false, $false, $false, $true, 1,
    $false, "Bubble Sort" + $cr, 2)

$extra = $d.Paragraphs.Item($algoritmoIndex + 1)
$extra.Range.Delete()

# ---------------------------------------------------------------------------
# 3. Table cell "CPU" -> "Tiempo de Ejecución"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "CPU", $true, $false, $false, $false, $false, $true, 1, $false,
    "Tiempo de Ejecución", 2)
